$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(29, 8).Value = 2450.6667
$ws.Cells.Item(29, 9).Value = 2450.6667
$ws.Cells.Item(29, 11).Value = 7352.000100000001
$ws.Cells.Item(29, 13).Value = -7071.000100000001
$ws.Cells.Item(86, 8).Value = 8085.2856
$ws.Cells.Item(86, 9).Value = 9961.875
$ws.Cells.Item(86, 10).Value = 5583.1665
$ws.Cells.Item(86, 11).Value = 9961.875
$ws.Cells.Item(86, 12).Value = 5583.1665
$ws.Cells.Item(86, 13).Value = -8838.875
$ws.Cells.Item(86, 14).Value = -7829.1665
$ws.Cells.Item(89, 8).Value = 8085.2856
$ws.Cells.Item(89, 9).Value = 9961.875
$ws.Cells.Item(89, 10).Value = 5583.1665
$ws.Cells.Item(89, 11).Value = 49809.375
$ws.Cells.Item(89, 12).Value = 27915.8325
$ws.Cells.Item(89, 13).Value = -44193.375
$ws.Cells.Item(89, 14).Value = -39147.8325
$ws.Cells.Item(138, 8).Value = 3684.3845
$ws.Cells.Item(138, 9).Value = 3254
$ws.Cells.Item(138, 10).Value = 3858.8647
$ws.Cells.Item(138, 11).Value = 9762
$ws.Cells.Item(138, 12).Value = 11576.5941
$ws.Cells.Item(138, 13).Value = -4622
$ws.Cells.Item(138, 14).Value = -21856.5941
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 8753.195
$ws.Cells.Item(32, 9).Value = 7092.9863
$ws.Cells.Item(32, 11).Value = 7092.9863
$ws.Cells.Item(32, 13).Value = -6805.9863
$ws.Cells.Item(74, 8).Value = 256972
$ws.Cells.Item(74, 10).Value = 309296
$ws.Cells.Item(74, 12).Value = 309296
$ws.Cells.Item(74, 14).Value = -311044
$ws.Cells.Item(77, 8).Value = 256972
$ws.Cells.Item(77, 10).Value = 309296
$ws.Cells.Item(77, 12).Value = 1546480
$ws.Cells.Item(77, 14).Value = -1555216
$ws.Cells.Item(88, 8).Value = 1838.9445
$ws.Cells.Item(88, 9).Value = 2736.8572
$ws.Cells.Item(88, 10).Value = 1267.5454
$ws.Cells.Item(88, 11).Value = 2736.8572
$ws.Cells.Item(88, 12).Value = 1267.5454
$ws.Cells.Item(88, 13).Value = -2330.8572
$ws.Cells.Item(88, 14).Value = -2079.5454
$ws.Cells.Item(91, 8).Value = 1838.9445
$ws.Cells.Item(91, 9).Value = 2736.8572
$ws.Cells.Item(91, 10).Value = 1267.5454
$ws.Cells.Item(91, 11).Value = 2736.8572
$ws.Cells.Item(91, 12).Value = 1267.5454
$ws.Cells.Item(91, 13).Value = -1332.8572
$ws.Cells.Item(91, 14).Value = -4075.5454
$ws.Cells.Item(97, 8).Value = 1598523.2
$ws.Cells.Item(97, 9).Value = 2097992.2
$ws.Cells.Item(97, 11).Value = 2097992.2
$ws.Cells.Item(97, 13).Value = -2097496.2
$ws.Cells.Item(132, 8).Value = 50894.477
$ws.Cells.Item(132, 9).Value = 18648.666
$ws.Cells.Item(132, 11).Value = 55945.99800000001
$ws.Cells.Item(132, 13).Value = -53415.99800000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 20643644
$ws.Cells.Item(86, 9).Value = 54169670
$ws.Cells.Item(86, 10).Value = 12242.77
$ws.Cells.Item(86, 11).Value = 54169670
$ws.Cells.Item(86, 12).Value = 12242.77
$ws.Cells.Item(86, 13).Value = -54168547
$ws.Cells.Item(86, 14).Value = -14488.77
$ws.Cells.Item(89, 8).Value = 20643644
$ws.Cells.Item(89, 9).Value = 54169670
$ws.Cells.Item(89, 10).Value = 12242.77
$ws.Cells.Item(89, 11).Value = 270848350
$ws.Cells.Item(89, 12).Value = 61213.85000000001
$ws.Cells.Item(89, 13).Value = -270842734
$ws.Cells.Item(89, 14).Value = -72445.85000000001
$ws.Cells.Item(105, 8).Value = 3474484.5
$ws.Cells.Item(105, 9).Value = 3908545.8
$ws.Cells.Item(105, 11).Value = 3908545.8
$ws.Cells.Item(105, 13).Value = -3906798.8
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(76, 8).Value = 5218.75
$ws.Cells.Item(76, 9).Value = 5218.75
$ws.Cells.Item(76, 11).Value = 5218.75
$ws.Cells.Item(76, 13).Value = -4903.75
$ws.Cells.Item(79, 8).Value = 5218.75
$ws.Cells.Item(79, 9).Value = 5218.75
$ws.Cells.Item(79, 11).Value = 5218.75
$ws.Cells.Item(79, 13).Value = -4126.75
$ws.Cells.Item(125, 8).Value = 39995
$ws.Cells.Item(125, 10).Value = 39995
$ws.Cells.Item(125, 12).Value = 39995
$ws.Cells.Item(125, 14).Value = -44915
$ws.Cells.Item(132, 8).Value = 72820.74000000001
$ws.Cells.Item(132, 9).Value = 48023.727
$ws.Cells.Item(132, 11).Value = 144071.181
$ws.Cells.Item(132, 13).Value = -141541.181
$ws.Cells.Item(141, 8).Value = 124900.2
$ws.Cells.Item(141, 10).Value = 131964.83
$ws.Cells.Item(141, 12).Value = 131964.83
$ws.Cells.Item(141, 14).Value = -142324.83
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 1357
$ws.Cells.Item(34, 10).Value = 2466.6667
$ws.Cells.Item(34, 12).Value = 7400.000100000001
$ws.Cells.Item(34, 14).Value = -7568.000100000001
$ws.Cells.Item(36, 8).Value = 765
$ws.Cells.Item(36, 9).Value = 765
$ws.Cells.Item(36, 11).Value = 2295
$ws.Cells.Item(36, 13).Value = -2126
$ws.Cells.Item(88, 8).Value = 16333.333
$ws.Cells.Item(88, 9).Value = 19000
$ws.Cells.Item(88, 10).Value = 15000
$ws.Cells.Item(88, 11).Value = 57000
$ws.Cells.Item(88, 12).Value = 45000
$ws.Cells.Item(88, 13).Value = -56572
$ws.Cells.Item(88, 14).Value = -45856
$ws.Cells.Item(91, 8).Value = 16333.333
$ws.Cells.Item(91, 9).Value = 19000
$ws.Cells.Item(91, 10).Value = 15000
$ws.Cells.Item(91, 11).Value = 57000
$ws.Cells.Item(91, 12).Value = 45000
$ws.Cells.Item(91, 13).Value = -55518
$ws.Cells.Item(91, 14).Value = -47964
$ws.Cells.Item(122, 8).Value = 1423.65
$ws.Cells.Item(122, 9).Value = 1313.3334
$ws.Cells.Item(122, 11).Value = 11820.0006
$ws.Cells.Item(122, 13).Value = -9370.000599999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(141, 8).Value = 109500.5
$ws.Cells.Item(141, 10).Value = 109500.5
$ws.Cells.Item(141, 12).Value = 109500.5
$ws.Cells.Item(141, 14).Value = -119860.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 11772
$ws.Cells.Item(40, 9).Value = 6126
$ws.Cells.Item(40, 11).Value = 6126
$ws.Cells.Item(40, 13).Value = -5990
$ws.Cells.Item(127, 8).Value = 80000
$ws.Cells.Item(127, 10).Value = 80000
$ws.Cells.Item(127, 12).Value = 80000
$ws.Cells.Item(127, 14).Value = -89920
$ws.Cells.Item(136, 8).Value = 99835.766
$ws.Cells.Item(136, 10).Value = 7998.5
$ws.Cells.Item(136, 12).Value = 23995.5
$ws.Cells.Item(136, 14).Value = -29095.5
$ws.Cells.Item(138, 8).Value = 82285.60000000001
$ws.Cells.Item(138, 10).Value = 82285.60000000001
$ws.Cells.Item(138, 12).Value = 82285.60000000001
$ws.Cells.Item(138, 14).Value = -92565.60000000001
$ws.Cells.Item(140, 8).Value = 99999
$ws.Cells.Item(140, 10).Value = 99999
$ws.Cells.Item(140, 12).Value = 99999
$ws.Cells.Item(140, 14).Value = -110359
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(110, 8).Value = 0
$ws.Cells.Item(110, 10).Value = 0
$ws.Cells.Item(110, 12).Value = 0
$ws.Cells.Item(110, 14).ClearContents()
$ws.Cells.Item(125, 8).Value = 90000
$ws.Cells.Item(125, 10).Value = 90000
$ws.Cells.Item(125, 12).Value = 90000
$ws.Cells.Item(125, 14).Value = -99840
$ws.Cells.Item(136, 8).Value = 5227.1284
$ws.Cells.Item(136, 9).Value = 5193.6313
$ws.Cells.Item(136, 11).Value = 15580.8939
$ws.Cells.Item(136, 13).Value = -13030.8939

Write-Host "done"